$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (THE PLAYERS Championship, Valspar Championship),
# shifting all rows below up by two.
$ws.Range("2:3").Delete()

$ws.Range("I10").Select()
